# Refresh crypto price + volume(1h) figures (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.786.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "'2.101.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'226.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").Value = "'62.20"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.46%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("D10").Value = "'0.0845"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "'15.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.89%  "
$ws.Range("D13").Value = "'2.415.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").Value = "'21.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "'0.800"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").Value = "'5.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "'2.125.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").Value = "'38.865.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "'71.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "'6.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").Value = "'227.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'2.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.02%  "
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("D26").Value = "'9.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.54%  "
$ws.Range("D27").Value = "'170.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("E29").Value = "  +3.29%  "
$ws.Range("D30").Value = "'19.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("E31").Value = "  +9.22%  "
$ws.Range("D32").Value = "'0.120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "'4.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").Value = "'7.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.95%  "
$ws.Range("D35").Value = "'4.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").Value = "'0.0613"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'3.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").Value = "'17.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("E41").Value = "  +3.46%  "
$ws.Range("D42").Value = "'101.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("D43").Value = "'1.525.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("E44").Value = "  +7.24%  "
$ws.Range("D45").Value = "'2.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("E48").Value = "  +5.15%  "
$ws.Range("D49").Value = "'4.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "'2.302.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.26%  "
